$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style for plain, unstyled text/number cells (no explicit cell style),
# used to strip the "quote prefix" style Excel applies when a numeric-looking
# string is entered, so the price cells keep their original default formatting.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "27.598.77"
$ws.Range("E2").Value = "  +3.10%  "
$ws.Range("D3").Value = "1.849.28"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("D4").Value = "'1.031"
$ws.Range("E4").Value = "  +2.86%  "
$ws.Range("D5").Value = "'320.82"
$ws.Range("E5").Value = "  +3.93%  "
$ws.Range("D6").Value = "'1.027"
$ws.Range("E6").Value = "  +2.54%  "
$ws.Range("D7").Value = "'0.4375"
$ws.Range("E7").Value = "  +1.67%  "
$ws.Range("D8").Value = "'0.3741"
$ws.Range("E8").Value = "  +2.02%  "
$ws.Range("D9").Value = "'0.07403"
$ws.Range("E9").Value = "  +3.01%  "
$ws.Range("D10").Value = "'0.8763"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").Value = "'21.47"
$ws.Range("E11").Value = "  +3.38%  "
$ws.Range("D12").Value = "1.865.68"
$ws.Range("E12").Value = "  -3.28%  "
$ws.Range("D13").Value = "'5.500"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("D14").Value = "'6.683"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "'0.07163"
$ws.Range("E15").Value = "  +4.14%  "
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "'1.031"
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "'0.000009026"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "'1.026"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").Value = "27.605.88"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("D22").Value = "'5.254"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "'11.21"
$ws.Range("D24").Value = "2.067.28"
$ws.Range("D25").Value = "'157.15"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").Value = "'1.932"
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").Value = "'5.292"
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").Value = "'1.952"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'0.09073"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("D32").Value = "'1.209"
$ws.Range("D33").Value = "'0.7679"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("D35").Value = "'2.877"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("D36").Value = "'1.028"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "'1.148"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").Value = "'0.01980"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "'0.05267"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.812"
$ws.Range("E40").Value = "  +6.63%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5171"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").Value = "'0.1673"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").Value = "'6.705"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("D44").Value = "'8.573"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").Value = "'108.97"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").Value = "'10.56"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "'1.716"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("E48").Value = "  +2.13%  "
$ws.Range("D49").Value = "'0.06379"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("D50").Value = "'1.889"
$ws.Range("E50").Value = "  +5.56%  "
$ws.Range("D51").Value = "'39.53"
$ws.Range("E51").Value = "  +6.11%  "

# Restore default (unstyled) formatting on price cells that were entered as
# quoted text so numeric-looking values (e.g. "1.031") are preserved verbatim.
$ws.Range("D4").Style = $defaultStyle
$ws.Range("D5").Style = $defaultStyle
$ws.Range("D6").Style = $defaultStyle
$ws.Range("D7").Style = $defaultStyle
$ws.Range("D8").Style = $defaultStyle
$ws.Range("D9").Style = $defaultStyle
$ws.Range("D10").Style = $defaultStyle
$ws.Range("D11").Style = $defaultStyle
$ws.Range("D13").Style = $defaultStyle
$ws.Range("D14").Style = $defaultStyle
$ws.Range("D15").Style = $defaultStyle
$ws.Range("D17").Style = $defaultStyle
$ws.Range("D18").Style = $defaultStyle
$ws.Range("D19").Style = $defaultStyle
$ws.Range("D22").Style = $defaultStyle
$ws.Range("D23").Style = $defaultStyle
$ws.Range("D25").Style = $defaultStyle
$ws.Range("D26").Style = $defaultStyle
$ws.Range("D28").Style = $defaultStyle
$ws.Range("D29").Style = $defaultStyle
$ws.Range("D31").Style = $defaultStyle
$ws.Range("D32").Style = $defaultStyle
$ws.Range("D33").Style = $defaultStyle
$ws.Range("D35").Style = $defaultStyle
$ws.Range("D36").Style = $defaultStyle
$ws.Range("D37").Style = $defaultStyle
$ws.Range("D38").Style = $defaultStyle
$ws.Range("D39").Style = $defaultStyle
$ws.Range("D40").Style = $defaultStyle
$ws.Range("D41").Style = $defaultStyle
$ws.Range("D42").Style = $defaultStyle
$ws.Range("D43").Style = $defaultStyle
$ws.Range("D44").Style = $defaultStyle
$ws.Range("D45").Style = $defaultStyle
$ws.Range("D46").Style = $defaultStyle
$ws.Range("D47").Style = $defaultStyle
$ws.Range("D49").Style = $defaultStyle
$ws.Range("D50").Style = $defaultStyle
$ws.Range("D51").Style = $defaultStyle
